$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "m2 x 7"
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = "metal"

$ws.Range("A1:D17").Select()
